$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Columns("T:W").Insert()
Write-Host "Sheet1 dims:" $ws.UsedRange.Address()
Write-Host "T1:" $ws.Range("T1").Value()
Write-Host "X1:" $ws.Range("X1").Value()
